# Apply per-cell value updates to Sheet1, matching the upstream
# coinranking.com price/volume refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.096.78"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "3.089.95"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.84"
$ws.Range("E5").Value = "  +2.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.24"
$ws.Range("E6").Value = "  +2.74%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.087.14"
$ws.Range("E8").Value = "  +0.99%  "
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("E11").Value = "  -5.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.470"
$ws.Range("E12").Value = "  +3.63%  "
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.07"
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("D15").Value = "3.594.19"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "64.162.91"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("D17").Value = "3.094.34"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.73"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("E21").Value = "  +1.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.673"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("E23").Value = "  +3.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.15"
$ws.Range("E24").Value = "  +12.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.13"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("E27").Value = "  +1.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.98"
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.06"
$ws.Range("E29").Value = "  +3.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.31"
$ws.Range("E31").Value = "  +0.71%  "
$ws.Range("E32").Value = "  -0.62%  "
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("E34").Value = "  -2.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "55.78"
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.20"
$ws.Range("E36").Value = "  +3.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "453.43"
$ws.Range("E37").Value = "  -3.15%  "
$ws.Range("E38").Value = "  +2.80%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.96"
$ws.Range("E39").Value = "  +15.34%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0820"
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").Value = "2.976.66"
$ws.Range("E41").Value = "  -3.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.24"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("E43").Value = "  -3.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.90"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("E45").Value = "  +2.25%  "
$ws.Range("E47").Value = "  +3.19%  "
$ws.Range("E48").Value = "  +1.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "119.92"
$ws.Range("E49").Value = "  +3.00%  "
$ws.Range("E50").Value = "  +1.09%  "
$ws.Range("E51").Value = "  +0.36%  "
